$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Ifnb1"
$ws.Cells.Item(2, 3).Value2 = "Ifnar2"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 1
$ws.Cells.Item(2, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 7).Value2 = 0.002067
$ws.Cells.Item(2, 8).Value2 = 0.006201
$ws.Cells.Item(2, 9).Value2 = 0.004794555500401285
$ws.Cells.Item(2, 10).Value2 = 0.004794555500401286
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 17.95291733333333
$ws.Cells.Item(2, 14).Value2 = 53.858752
$ws.Cells.Item(2, 15).Value2 = 0.07393528283891208
$ws.Cells.Item(2, 16).Value2 = 0.07489467518434184
$ws.Cells.Item(2, 17).Value2 = 0.03710868012799999
$ws.Cells.Item(2, 18).Value2 = 0.333978121152
$ws.Cells.Item(2, 19).Value2 = 0.0003544868170090307
$ws.Cells.Item(2, 20).Value2 = 0.0003590866768558539

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Ifnb1"
$ws.Cells.Item(3, 3).Value2 = "Ifnar2"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(3, 7).Value2 = 0.002067
$ws.Cells.Item(3, 8).Value2 = 0.006201
$ws.Cells.Item(3, 9).Value2 = 0.004794555500401285
$ws.Cells.Item(3, 10).Value2 = 0.004794555500401286
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 30.84182166666666
$ws.Cells.Item(3, 14).Value2 = 92.525465
$ws.Cells.Item(3, 15).Value2 = 0.1270155020409099
$ws.Cells.Item(3, 16).Value2 = 0.1286636691369156
$ws.Cells.Item(3, 17).Value2 = 0.063750045385
$ws.Cells.Item(3, 18).Value2 = 0.573750408465
$ws.Cells.Item(3, 19).Value2 = 0.0006089828739464752
$ws.Cells.Item(3, 20).Value2 = 0.0006168851025622098

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Ifnb1"
$ws.Cells.Item(4, 3).Value2 = "Ifnar2"
$ws.Cells.Item(4, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value2 = 1
$ws.Cells.Item(4, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 7).Value2 = 0.002067
$ws.Cells.Item(4, 8).Value2 = 0.006201
$ws.Cells.Item(4, 9).Value2 = 0.004794555500401285
$ws.Cells.Item(4, 10).Value2 = 0.004794555500401286
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 100.6030296666667
$ws.Cells.Item(4, 14).Value2 = 301.809089
$ws.Cells.Item(4, 15).Value2 = 0.4143122432277931
$ws.Cells.Item(4, 16).Value2 = 0.4196884043718117
$ws.Cells.Item(4, 17).Value2 = 0.207946462321
$ws.Cells.Item(4, 18).Value2 = 1.871518160889
$ws.Cells.Item(4, 19).Value2 = 0.001986443044651411
$ws.Cells.Item(4, 20).Value2 = 0.002012219347635509

# Row 5
$ws.Cells.Item(5, 1).Value2 = "ECs"
$ws.Cells.Item(5, 2).Value2 = "Ifnb1"
$ws.Cells.Item(5, 3).Value2 = "Ifnar2"
$ws.Cells.Item(5, 4).Value2 = "MuSCs"
$ws.Cells.Item(5, 5).Value2 = 1
$ws.Cells.Item(5, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(5, 7).Value2 = 0.002067
$ws.Cells.Item(5, 8).Value2 = 0.006201
$ws.Cells.Item(5, 9).Value2 = 0.004794555500401285
$ws.Cells.Item(5, 10).Value2 = 0.004794555500401286
$ws.Cells.Item(5, 11).Value2 = 2
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 9.3314655
$ws.Cells.Item(5, 14).Value2 = 18.662931
$ws.Cells.Item(5, 15).Value2 = 0.03842966177775806
$ws.Cells.Item(5, 16).Value2 = 0.02595221952474473
$ws.Cells.Item(5, 17).Value2 = 0.0192881391885
$ws.Cells.Item(5, 18).Value2 = 0.115728835131
$ws.Cells.Item(5, 19).Value2 = 0.000184253146255111
$ws.Cells.Item(5, 20).Value2 = 0.0001244293568699865

# Row 6
$ws.Cells.Item(6, 1).Value2 = "ECs"
$ws.Cells.Item(6, 2).Value2 = "Ifnb1"
$ws.Cells.Item(6, 3).Value2 = "Ifnar2"
$ws.Cells.Item(6, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value2 = 1
$ws.Cells.Item(6, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(6, 7).Value2 = 0.002067
$ws.Cells.Item(6, 8).Value2 = 0.006201
$ws.Cells.Item(6, 9).Value2 = 0.004794555500401285
$ws.Cells.Item(6, 10).Value2 = 0.004794555500401286
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 84.09011600000001
$ws.Cells.Item(6, 14).Value2 = 252.270348
$ws.Cells.Item(6, 15).Value2 = 0.3463073101146269
$ws.Cells.Item(6, 16).Value2 = 0.3508010317821862
$ws.Cells.Item(6, 17).Value2 = 0.173814269772
$ws.Cells.Item(6, 18).Value2 = 1.564328427948
$ws.Cells.Item(6, 19).Value2 = 0.001660389618539258
$ws.Cells.Item(6, 20).Value2 = 0.001681935016477727

# Row 7
$ws.Cells.Item(7, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value2 = "Ifnb1"
$ws.Cells.Item(7, 3).Value2 = "Ifnar2"
$ws.Cells.Item(7, 4).Value2 = "ECs"
$ws.Cells.Item(7, 5).Value2 = 2
$ws.Cells.Item(7, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(7, 7).Value2 = 0.3352766666666667
$ws.Cells.Item(7, 8).Value2 = 1.00583
$ws.Cells.Item(7, 9).Value2 = 0.7776983968664127
$ws.Cells.Item(7, 10).Value2 = 0.7776983968664127
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 17.95291733333333
$ws.Cells.Item(7, 14).Value2 = 53.858752
$ws.Cells.Item(7, 15).Value2 = 0.07393528283891208
$ws.Cells.Item(7, 16).Value2 = 0.07489467518434184
$ws.Cells.Item(7, 17).Value2 = 6.019194280462222
$ws.Cells.Item(7, 18).Value2 = 54.17274852415999
$ws.Cells.Item(7, 19).Value2 = 0.05749935093568673
$ws.Cells.Item(7, 20).Value2 = 0.05824546882469336

# Row 8
$ws.Cells.Item(8, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value2 = "Ifnb1"
$ws.Cells.Item(8, 3).Value2 = "Ifnar2"
$ws.Cells.Item(8, 4).Value2 = "FAPs"
$ws.Cells.Item(8, 5).Value2 = 2
$ws.Cells.Item(8, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(8, 7).Value2 = 0.3352766666666667
$ws.Cells.Item(8, 8).Value2 = 1.00583
$ws.Cells.Item(8, 9).Value2 = 0.7776983968664127
$ws.Cells.Item(8, 10).Value2 = 0.7776983968664127
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 30.84182166666666
$ws.Cells.Item(8, 14).Value2 = 92.525465
$ws.Cells.Item(8, 15).Value2 = 0.1270155020409099
$ws.Cells.Item(8, 16).Value2 = 0.1286636691369156
$ws.Cells.Item(8, 17).Value2 = 10.34054316232778
$ws.Cells.Item(8, 18).Value2 = 93.06488846095
$ws.Cells.Item(8, 19).Value2 = 0.0987797523143982
$ws.Cells.Item(8, 20).Value2 = 0.1000615292227298

# Row 9
$ws.Cells.Item(9, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value2 = "Ifnb1"
$ws.Cells.Item(9, 3).Value2 = "Ifnar2"
$ws.Cells.Item(9, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value2 = 2
$ws.Cells.Item(9, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(9, 7).Value2 = 0.3352766666666667
$ws.Cells.Item(9, 8).Value2 = 1.00583
$ws.Cells.Item(9, 9).Value2 = 0.7776983968664127
$ws.Cells.Item(9, 10).Value2 = 0.7776983968664127
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 100.6030296666667
$ws.Cells.Item(9, 14).Value2 = 301.809089
$ws.Cells.Item(9, 15).Value2 = 0.4143122432277931
$ws.Cells.Item(9, 16).Value2 = 0.4196884043718117
$ws.Cells.Item(9, 17).Value2 = 33.72984844320777
$ws.Cells.Item(9, 18).Value2 = 303.56863598887
$ws.Cells.Item(9, 19).Value2 = 0.3222099673603819
$ws.Cells.Item(9, 20).Value2 = 0.3263909992633807

# Row 10
$ws.Cells.Item(10, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value2 = "Ifnb1"
$ws.Cells.Item(10, 3).Value2 = "Ifnar2"
$ws.Cells.Item(10, 4).Value2 = "MuSCs"
$ws.Cells.Item(10, 5).Value2 = 2
$ws.Cells.Item(10, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(10, 7).Value2 = 0.3352766666666667
$ws.Cells.Item(10, 8).Value2 = 1.00583
$ws.Cells.Item(10, 9).Value2 = 0.7776983968664127
$ws.Cells.Item(10, 10).Value2 = 0.7776983968664127
$ws.Cells.Item(10, 11).Value2 = 2
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 9.3314655
$ws.Cells.Item(10, 14).Value2 = 18.662931
$ws.Cells.Item(10, 15).Value2 = 0.03842966177775806
$ws.Cells.Item(10, 16).Value2 = 0.02595221952474473
$ws.Cells.Item(10, 17).Value2 = 3.128622647955
$ws.Cells.Item(10, 18).Value2 = 18.77173588773
$ws.Cells.Item(10, 19).Value2 = 0.0298866863566809
$ws.Cells.Item(10, 20).Value2 = 0.02018299951951919

# Row 11
$ws.Cells.Item(11, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value2 = "Ifnb1"
$ws.Cells.Item(11, 3).Value2 = "Ifnar2"
$ws.Cells.Item(11, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value2 = 2
$ws.Cells.Item(11, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(11, 7).Value2 = 0.3352766666666667
$ws.Cells.Item(11, 8).Value2 = 1.00583
$ws.Cells.Item(11, 9).Value2 = 0.7776983968664127
$ws.Cells.Item(11, 10).Value2 = 0.7776983968664127
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 84.09011600000001
$ws.Cells.Item(11, 14).Value2 = 252.270348
$ws.Cells.Item(11, 15).Value2 = 0.3463073101146269
$ws.Cells.Item(11, 16).Value2 = 0.3508010317821862
$ws.Cells.Item(11, 17).Value2 = 28.19345379209334
$ws.Cells.Item(11, 18).Value2 = 253.74108412884
$ws.Cells.Item(11, 19).Value2 = 0.269322639899265
$ws.Cells.Item(11, 20).Value2 = 0.2728174000360897

# Row 12
$ws.Cells.Item(12, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value2 = "Ifnb1"
$ws.Cells.Item(12, 3).Value2 = "Ifnar2"
$ws.Cells.Item(12, 4).Value2 = "ECs"
$ws.Cells.Item(12, 5).Value2 = 1
$ws.Cells.Item(12, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(12, 7).Value2 = 0.09377033333333333
$ws.Cells.Item(12, 8).Value2 = 0.281311
$ws.Cells.Item(12, 9).Value2 = 0.217507047633186
$ws.Cells.Item(12, 10).Value2 = 0.217507047633186
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 12).Value2 = 1
$ws.Cells.Item(12, 13).Value2 = 17.95291733333333
$ws.Cells.Item(12, 14).Value2 = 53.858752
$ws.Cells.Item(12, 15).Value2 = 0.07393528283891208
$ws.Cells.Item(12, 16).Value2 = 0.07489467518434184
$ws.Cells.Item(12, 17).Value2 = 1.683451042652444
$ws.Cells.Item(12, 18).Value2 = 15.151059383872
$ws.Cells.Item(12, 19).Value2 = 0.01608144508621633
$ws.Cells.Item(12, 20).Value2 = 0.01629011968279263

# Row 13
$ws.Cells.Item(13, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value2 = "Ifnb1"
$ws.Cells.Item(13, 3).Value2 = "Ifnar2"
$ws.Cells.Item(13, 4).Value2 = "FAPs"
$ws.Cells.Item(13, 5).Value2 = 1
$ws.Cells.Item(13, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(13, 7).Value2 = 0.09377033333333333
$ws.Cells.Item(13, 8).Value2 = 0.281311
$ws.Cells.Item(13, 9).Value2 = 0.217507047633186
$ws.Cells.Item(13, 10).Value2 = 0.217507047633186
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 12).Value2 = 1
$ws.Cells.Item(13, 13).Value2 = 30.84182166666666
$ws.Cells.Item(13, 14).Value2 = 92.525465
$ws.Cells.Item(13, 15).Value2 = 0.1270155020409099
$ws.Cells.Item(13, 16).Value2 = 0.1286636691369156
$ws.Cells.Item(13, 17).Value2 = 2.892047898290555
$ws.Cells.Item(13, 18).Value2 = 26.028431084615
$ws.Cells.Item(13, 19).Value2 = 0.02762676685256522
$ws.Cells.Item(13, 20).Value2 = 0.02798525481162357

# Row 14
$ws.Cells.Item(14, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value2 = "Ifnb1"
$ws.Cells.Item(14, 3).Value2 = "Ifnar2"
$ws.Cells.Item(14, 4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value2 = 1
$ws.Cells.Item(14, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(14, 7).Value2 = 0.09377033333333333
$ws.Cells.Item(14, 8).Value2 = 0.281311
$ws.Cells.Item(14, 9).Value2 = 0.217507047633186
$ws.Cells.Item(14, 10).Value2 = 0.217507047633186
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 12).Value2 = 1
$ws.Cells.Item(14, 13).Value2 = 100.6030296666667
$ws.Cells.Item(14, 14).Value2 = 301.809089
$ws.Cells.Item(14, 15).Value2 = 0.4143122432277931
$ws.Cells.Item(14, 16).Value2 = 0.4196884043718117
$ws.Cells.Item(14, 17).Value2 = 9.433579626186555
$ws.Cells.Item(14, 18).Value2 = 84.90221663567898
$ws.Cells.Item(14, 19).Value2 = 0.09011583282275971
$ws.Cells.Item(14, 20).Value2 = 0.09128518576079545

# Row 15
$ws.Cells.Item(15, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value2 = "Ifnb1"
$ws.Cells.Item(15, 3).Value2 = "Ifnar2"
$ws.Cells.Item(15, 4).Value2 = "MuSCs"
$ws.Cells.Item(15, 5).Value2 = 1
$ws.Cells.Item(15, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(15, 7).Value2 = 0.09377033333333333
$ws.Cells.Item(15, 8).Value2 = 0.281311
$ws.Cells.Item(15, 9).Value2 = 0.217507047633186
$ws.Cells.Item(15, 10).Value2 = 0.217507047633186
$ws.Cells.Item(15, 11).Value2 = 2
$ws.Cells.Item(15, 12).Value2 = 1
$ws.Cells.Item(15, 13).Value2 = 9.3314655
$ws.Cells.Item(15, 14).Value2 = 18.662931
$ws.Cells.Item(15, 15).Value2 = 0.03842966177775806
$ws.Cells.Item(15, 16).Value2 = 0.02595221952474473
$ws.Cells.Item(15, 17).Value2 = 0.8750146304235
$ws.Cells.Item(15, 18).Value2 = 5.250087782541
$ws.Cells.Item(15, 19).Value2 = 0.00835872227482205
$ws.Cells.Item(15, 20).Value2 = 0.00564479064835555

# Row 16
$ws.Cells.Item(16, 1).Value2 = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value2 = "Ifnb1"
$ws.Cells.Item(16, 3).Value2 = "Ifnar2"
$ws.Cells.Item(16, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value2 = 1
$ws.Cells.Item(16, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(16, 7).Value2 = 0.09377033333333333
$ws.Cells.Item(16, 8).Value2 = 0.281311
$ws.Cells.Item(16, 9).Value2 = 0.217507047633186
$ws.Cells.Item(16, 10).Value2 = 0.217507047633186
$ws.Cells.Item(16, 11).Value2 = 3
$ws.Cells.Item(16, 12).Value2 = 1
$ws.Cells.Item(16, 13).Value2 = 84.09011600000001
$ws.Cells.Item(16, 14).Value2 = 252.270348
$ws.Cells.Item(16, 15).Value2 = 0.3463073101146269
$ws.Cells.Item(16, 16).Value2 = 0.3508010317821862
$ws.Cells.Item(16, 17).Value2 = 7.885158207358668
$ws.Cells.Item(16, 18).Value2 = 70.966423866228
$ws.Cells.Item(16, 19).Value2 = 0.07532428059682265
$ws.Cells.Item(16, 20).Value2 = 0.07630169672961876
